$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Refresh the cached "datetimeFigureOut" footer date shown on the slide
#    master and every slide layout (2021-08-05 -> 2021-08-08).
# ---------------------------------------------------------------------------
function Update-DatePlaceholders($shapes) {
  for ($i = 1; $i -le $shapes.Count; $i++) {
    $sh = $shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
      if ($sh.TextFrame.TextRange.Text -eq "2021-08-05") {
        $sh.TextFrame.TextRange.Text = "2021-08-08"
      }
    }
  }
}

Update-DatePlaceholders $p.SlideMaster.Shapes

$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
  Update-DatePlaceholders $layouts.Item($L).Shapes
}

# ---------------------------------------------------------------------------
# 2) Slide 17 ("Object List" slide): drop four of the rectangle callouts
#    (슈퍼 / 서브 8 / 스케일 / 테라인) and slide the remaining "오브젝트 리스트"
#    rectangle up into the now-empty top-middle slot.
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(17)

foreach ($name in @("직사각형 3", "직사각형 4", "직사각형 6", "직사각형 7")) {
  $s.Shapes.Item($name).Delete()
}

$moved = $s.Shapes.Item("직사각형 8")
$moved.Left = 316.3956762913386   # -> 4018225 EMU
$moved.Top  = 158.400002          # -> 2011680 EMU
